$d = $word.ActiveDocument

# Locate the list-item paragraph containing
#   "Comentarios sobre las elecciones y dificultades "
# (the numbered-list occurrence, not the earlier index/TOC mention of the
# same words).
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Comentarios sobre las elecciones y dificultades " -and `
        $p.Range.ParagraphFormat.Style.NameLocal -eq "List Paragraph") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph"
}

# Merge this paragraph with the following one ("encontradas en el
# desarrollo del proyecto") by deleting the paragraph mark that separates
# them. This keeps the first paragraph's pPr/numbering and produces a
# single paragraph whose text is the concatenation of both runs.
$end = $target.Range.End
$markRange = $d.Range($end - 1, $end)
$markRange.Delete()

# The merged paragraph now reads:
#   "Comentarios sobre las elecciones y dificultades encontradas en el desarrollo del proyecto"
# Re-author its run content as three runs - the heading text (no trailing
# space), a standalone space run, and the existing continuation text -
# using InsertXML so the <w:lastRenderedPageBreak/> marker that belongs on
# the first run is preserved explicitly (simple Range-splitting APIs drop
# it when a run gets divided).
$paraStart = $target.Range.Start
$paraEnd = $target.Range.End - 1   # exclude the paragraph mark
$contentRange = $d.Range($paraStart, $paraEnd)

$openXmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$runProps = '<w:rPr><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr>'
$run1 = '<w:r w:rsidRPr="00F93B14">' + $runProps + '<w:lastRenderedPageBreak/><w:t>Comentarios sobre las elecciones y dificultades</w:t></w:r>'
$run2 = '<w:r w:rsidRPr="00F93B14">' + $runProps + '<w:t xml:space="preserve"> </w:t></w:r>'
$run3 = '<w:r w:rsidRPr="00F93B14">' + $runProps + '<w:t>encontradas en el desarrollo del proyecto</w:t></w:r>'

$fragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
    '<w:document ' + $openXmlNs + '><w:body><w:p>' + $run1 + $run2 + $run3 + '</w:p></w:body></w:document>' + `
    '</pkg:xmlData></pkg:part></pkg:package>'

$contentRange.InsertXML($fragment)
